$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update District (column G) values to official names from website
$ws.Range("G3").Value  = "Kalaburagi (Gulbarga)"
$ws.Range("G5").Value  = "Kalaburagi (Gulbarga)"
$ws.Range("G6").Value  = "Kalaburagi (Gulbarga)"
$ws.Range("G19").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G25").Value = "Ramanagara"
$ws.Range("G27").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G28").Value = "Ballari (Bellary)"
$ws.Range("G29").Value = "Davangere"
$ws.Range("G31").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G35").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G36").Value = "Vijayapura (Bijapur)"
$ws.Range("G41").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G43").Value = "Ramanagara"
$ws.Range("G48").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G56").Value = "Kalaburagi (Gulbarga)"

# Remove the stray empty cells in column F (rows 34 and 52) that had no content
$ws.Range("F34").ClearContents()
$ws.Range("F52").ClearContents()
